$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Donor cells that already carry the exact style used by the target cells and
# are never themselves modified by this edit - used to restore the original
# cell style after a numeric-looking string has been entered (Excel's COM
# layer auto-converts a bare numeric-looking string like "1.0" into a real
# number unless the cell is Text-formatted first; re-pasting the original
# formatting afterwards keeps the display/style identical while the stored
# value stays the literal text "1.0").
$donorStyle10 = $ws.Range("C21")   # style s="10" (At Work / Sick Leave / Childcare / Annual columns)
$donorStyle16 = $ws.Range("D44")   # style s="16" (Total row)

function Set-TextValue($rangeRef, $text, $donor) {
    $c = $ws.Range($rangeRef)
    $c.NumberFormat = "@"
    $c.Value = $text
    $donor.Copy()
    $c.PasteSpecial(-4122)   # xlPasteFormats - restore original look/style only
}

function Clear-Cell($rangeRef) {
    $ws.Range($rangeRef).Value = ""
}

# --- Daily rows: move the "1.0" mark between columns -------------------

# 06 - 09 January: move mark from "At Work" (C) to "Sick Leave" (E)
Clear-Cell "C16"
Set-TextValue "E16" "1.0" $donorStyle10

Clear-Cell "C17"
Set-TextValue "E17" "1.0" $donorStyle10

Clear-Cell "C18"
Set-TextValue "E18" "1.0" $donorStyle10

Clear-Cell "C19"
Set-TextValue "E19" "1.0" $donorStyle10

# 10 January: move mark from "Childcare Leave" (F) to "Sick Leave" (E)
Set-TextValue "E20" "1.0" $donorStyle10
Clear-Cell "F20"

# 13 - 15 January: move mark from "Childcare Leave" (F) to "At Work" (C)
Set-TextValue "C23" "1.0" $donorStyle10
Clear-Cell "F23"

Set-TextValue "C24" "1.0" $donorStyle10
Clear-Cell "F24"

Set-TextValue "C25" "1.0" $donorStyle10
Clear-Cell "F25"

# 16 - 17 January: move mark from "Childcare Leave" (F) and "Annual Leave" (G)
# to "At Work" (C)
Set-TextValue "C26" "1.0" $donorStyle10
Clear-Cell "F26"
Clear-Cell "G26"

Set-TextValue "C27" "1.0" $donorStyle10
Clear-Cell "F27"
Clear-Cell "G27"

# 20 - 23 January: move mark from "Annual Leave" (G) to "At Work" (C)
Set-TextValue "C30" "1.0" $donorStyle10
Clear-Cell "G30"

Set-TextValue "C31" "1.0" $donorStyle10
Clear-Cell "G31"

Set-TextValue "C32" "1.0" $donorStyle10
Clear-Cell "G32"

Set-TextValue "C33" "1.0" $donorStyle10
Clear-Cell "G33"

# 24 January: move mark from "Annual Leave" (G) to "Childcare Leave" (F)
Set-TextValue "F34" "1.0" $donorStyle10
Clear-Cell "G34"

# 27 January: move mark from "Annual Leave" (G) to "Childcare Leave" (F)
Set-TextValue "F37" "1.0" $donorStyle10
Clear-Cell "G37"

# 28 January: move mark from "Annual Leave" (G) to "At Work" (C)
Set-TextValue "C38" "1.0" $donorStyle10
Clear-Cell "G38"

# 31 January: move mark from "Annual Leave" (G) to "At Work" (C)
Set-TextValue "C41" "1.0" $donorStyle10
Clear-Cell "G41"

# --- Totals row ----------------------------------------------------------
Set-TextValue "C44" "11.0" $donorStyle16
Set-TextValue "E44" "7.0" $donorStyle16
Set-TextValue "F44" "2.0" $donorStyle16
Set-TextValue "G44" "-" $donorStyle16

$excel.CutCopyMode = $false
